$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their original text formatting
# (values such as "324.40" or "0.08160" must remain text, not be coerced to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "324.40"
$ws.Range("E2").Value = "9.23%"
$ws.Range("D3").Value = "49.80"
$ws.Range("E3").Value = "19.06%"
$ws.Range("D4").Value = "5.342"
$ws.Range("E4").Value = "6.84%"
$ws.Range("D5").Value = "0.08160"
$ws.Range("E5").Value = "8.54%"
$ws.Range("D6").Value = "4.613"
$ws.Range("E6").Value = "5.45%"
$ws.Range("D7").Value = "1.666"
$ws.Range("E7").Value = "5.23%"
$ws.Range("D8").Value = "1.174"
$ws.Range("E8").Value = "26.85%"
$ws.Range("D9").Value = "0.1343"
$ws.Range("E9").Value = "12.49%"
$ws.Range("D10").Value = "0.1958"
$ws.Range("E10").Value = "7.57%"
$ws.Range("D11").Value = "0.09499"
$ws.Range("E11").Value = "6.68%"
$ws.Range("D12").Value = "0.04552"
$ws.Range("E12").Value = "11.70%"
$ws.Range("D13").Value = "0.1047"
$ws.Range("E13").Value = "-0.19%"
$ws.Range("D14").Value = "0.001330"
$ws.Range("E14").Value = "3.96%"
$ws.Range("D15").Value = "0.005971"
$ws.Range("E15").Value = "2.34%"
$ws.Range("D16").Value = "3.398"
$ws.Range("E16").Value = "1.30%"
$ws.Range("D17").Value = "2.439"
$ws.Range("E17").Value = "1.57%"
$ws.Range("D18").Value = "0.3395"
$ws.Range("E18").Value = "2.46%"
$ws.Range("D19").Value = "8.181"
$ws.Range("E19").Value = "1.07%"
$ws.Range("E20").Value = "2.07%"
$ws.Range("D21").Value = "0.3052"
$ws.Range("E21").Value = "-1.62%"
$ws.Range("D22").Value = "0.04307"
$ws.Range("E22").Value = "4.94%"
$ws.Range("D23").Value = "0.001306"
$ws.Range("E23").Value = "3.01%"
$ws.Range("D24").Value = "0.004267"
$ws.Range("E24").Value = "9.45%"
$ws.Range("D25").Value = "0.0001348"
$ws.Range("E25").Value = "9.52%"
$ws.Range("D26").Value = "0.0003721"
$ws.Range("E26").Value = "-0.09%"
$ws.Range("D38").Value = "0.02779"
$ws.Range("E38").Value = "15.50%"
$ws.Range("D39").Value = "0.05557"
$ws.Range("E39").Value = "6.84%"
$ws.Range("D40").Value = "0.006293"
$ws.Range("E40").Value = "-0.21%"
$ws.Range("D41").Value = "0.007688"
$ws.Range("E41").Value = "-1.42%"
$ws.Range("D42").Value = "0.1447"
$ws.Range("E42").Value = "9.14%"
$ws.Range("D43").Value = "0.007691"
$ws.Range("E43").Value = "3.80%"
$ws.Range("D44").Value = "0.008065"
$ws.Range("E44").Value = "10.93%"
$ws.Range("D45").Value = "0.3520"
$ws.Range("E45").Value = "18.72%"
$ws.Range("D46").Value = "0.00006774"
$ws.Range("E46").Value = "2.70%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.21%"
$ws.Range("E48").Value = "93.87%"
$ws.Range("D49").Value = "0.003999"
$ws.Range("E49").Value = "-4.86%"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "-0.21%"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "-0.21%"
